$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: A3 changes text to "user@example.com" and gets a hyperlink
$ws.Range("A3").Formula = "user@example.com"
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:user@example.com")
# restore exact style (Hyperlink + left aligned, s4) by copying from A4
$ws.Range("A4").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# Row 4: A4 changes text to "customer1@example.com", keeps its existing style (s4, hyperlink+left)
$ws.Range("A4").Formula = "customer1@example.com"

# Row 5 (new): A5 = hatest@example.com, with hyperlink, style s1 (Hyperlink only)
$ws.Range("A5").Formula = "hatest@example.com"
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:hatest@example.com")
$ws.Range("A5").Style = "Hyperlink"
$ws.Range("B5").Formula = 123456

# Row 6 (new): A6 = hatest2@example.com, with hyperlink, style s1 (Hyperlink only)
$ws.Range("A6").Formula = "hatest2@example.com"
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:hatest2@example.com")
$ws.Range("A6").Style = "Hyperlink"
$ws.Range("B6").Formula = 123456

# Row 7 (new): A7 = admin@example.com, style s3 (copy from A2)
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Formula = "admin@example.com"
$ws.Range("B7").Formula = 123

# Row 8 (new): A8 = admin@example.com, style s3
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Formula = "admin@example.com"
$ws.Range("B8").Formula = 3434

# Update selection to B8
$ws.Range("B8").Select() | Out-Null

Write-Host "done"
